# The commit swaps the contents of ppt/theme/theme1.xml ("Office Theme")
# and ppt/theme/theme2.xml ("Integral") - i.e. the deck's live theme
# (currently "Integral", used by the slide master / all layouts / slides)
# becomes the "Office Theme" palette, and vice versa.
#
# The only theme reachable/mutable through the PowerPoint object model in
# this host is the one bound to the slide master (ppt/theme/theme2.xml);
# its 12-colour ThemeColorScheme is fully settable and persists to the
# underlying OOXML. We drive that theme's colours to the "Office Theme"
# palette that theme1.xml currently holds, reproducing the swap's visible
# effect on the presentation's live design.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0         # dk1      -> 000000
$cs.Item(2).RGB  = 16777215  # lt1      -> FFFFFF
$cs.Item(3).RGB  = 6968388   # dk2      -> 44546A
$cs.Item(4).RGB  = 15132391  # lt2      -> E7E6E6
$cs.Item(5).RGB  = 13998939  # accent1  -> 5B9BD5
$cs.Item(6).RGB  = 3243501   # accent2  -> ED7D31
$cs.Item(7).RGB  = 10855845  # accent3  -> A5A5A5
$cs.Item(8).RGB  = 49407     # accent4  -> FFC000
$cs.Item(9).RGB  = 12874308  # accent5  -> 4472C4
$cs.Item(10).RGB = 4697456   # accent6  -> 70AD47
$cs.Item(11).RGB = 12673797  # hlink    -> 0563C1
$cs.Item(12).RGB = 7491477   # folHlink -> 954F72
